# Apply the edit: replace the "Scherm 4/Baron...", "Foto van het hoofdgebouw...", and
# "Scherm 5/" paragraphs (the tail of the document) with the expanded set of
# Scherm 4-14 paragraphs described in the commit.

$d = $word.ActiveDocument

# Locate the three paragraphs that form the tail of the document, starting at
# "Scherm 4/Baron" and ending with the trailing "Scherm 5/" paragraph (which
# holds the _GoBack bookmark).
$total = $d.Paragraphs.Count
$firstPara = $d.Paragraphs($total - 2)
$lastPara = $d.Paragraphs($total)

if ($firstPara.Range.Text -notlike "Scherm 4/Baron*") {
    throw "Unexpected document shape: first paragraph was [$($firstPara.Range.Text)]"
}

$delRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)
$delRange.Delete()

$insPos = $d.Content.End
$insRange = $d.Range($insPos, $insPos)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scherm 4/Joris en de Draak/gewonnen-verloren</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>Foto van de gehele attractie, van de ingang richting de optakeling (..)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Scherm </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>/Baron</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>1898/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>gaJeMee</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>Foto van het hoofdgebouw met de ingang (</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>..)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Scherm </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>6</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Baron1898/mijn</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Droomvlucht/Fata </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Morgana</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> een foto van een van de grotten(…)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scherm 7/Baron1898/overleeft:</w:t></w:r></w:p><w:p><w:r><w:t>Foto van de medewerkers/rolstoel ingang</w:t></w:r><w:r><w:t xml:space="preserve"> (…)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scherm 8/Baron898/dood</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>Foto van voorshow 2 deuren (...)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scherm9/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>VogelRok</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>/wachtrij</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>Foto van wachtrij met daarop een bord met 0 min</w:t></w:r><w:r><w:t xml:space="preserve"> (…)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scherm10/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>VogelRok</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>/dood</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>Een medewerker met een mes</w:t></w:r><w:r><w:t xml:space="preserve"> (…)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scherm11/</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Vogelrok/gewonnen</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Foto van de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>unox</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> kraam in reizenrijk</w:t></w:r><w:r><w:t xml:space="preserve"> (…)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scherm12/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Pirana</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>/begin</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Foto van ingang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pirana</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (…)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scherm13/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Pirana</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>/dood</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>-overleven</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t xml:space="preserve">Een bewerkte foto van het </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pirana</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> eiland</w:t></w:r><w:r><w:t xml:space="preserve"> (…)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scherm14/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Symbolica</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>ontgrendeld</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>Foto van dichte voorshow deur</w:t></w:r><w:r><w:t xml:space="preserve"> (…)</w:t></w:r></w:p><w:p><w:r><w:t>Foto van koningszaal</w:t></w:r><w:r><w:t xml:space="preserve"> (…)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insRange.InsertXML($xml)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
